$wb = $excel.ActiveWorkbook

# --- Sheet "Reguły": reorder the object lists inside each rule description ---
$rules = $wb.Worksheets.Item("Reguły")
$rules.Range("B2").Value = "(attempts >=  3.0) & (pregnancy <=  0.0) => (class <= 1) ['a3', 'a7', 'a1']"
$rules.Range("B3").Value = "(sperm >=  3.0) => (class <= 1) ['a22', 'a25']"
$rules.Range("B4").Value = "(age >=  40.0) & (pregnancy <=  0.0) => (class <= 1) ['a3', 'a15']"
$rules.Range("B6").Value = "(age >=  42.0) => (class <= 1) ['a14', 'a3']"
$rules.Range("B7").Value = "(age <=  31.0) & (attempts <=  1.0) & (endometrium <=  1.0) => (class >= 2) ['a24', 'a9', 'a11', 'a12']"
$rules.Range("B8").Value = "(frozen_embryos >=  8.0) & (sperm <=  1.0) => (class >= 2) ['a16', 'a6']"

# --- Sheet "Statystyki reguł": update coverage (column C) values ---
$stats = $wb.Worksheets.Item("Statystyki reguł")
$stats.Range("C2").Value = 0.375
$stats.Range("C4").Value = 0.25

# --- Sheet "Walidacja krzyżowa": re-run produced rows in a new order ---
$cv = $wb.Worksheets.Item("Walidacja krzyżowa")
$cv.Range("A1").Value = "accuracy"
$cv.Range("B1").Value = 0.36
$cv.Range("A2").Value = "not_classified"
$cv.Range("B2").Value = 0.5600000000000001
$cv.Range("A3").Value = "correct"
$cv.Range("B3").Value = 0.8181818181818182
$cv.Range("A4").Value = "f1_score"
$cv.Range("B4").Value = 0.48
